$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure data range is treated as text so values like "$22.99" and "7.55%"
# are stored as literal strings rather than being coerced to numbers.
$ws.Range("B2:P10").NumberFormat = "@"

# Row 2: BIHLERFLEX
$ws.Range("A2").Value = "BIHLERFLEX"
$ws.Range("B2").Value = "$22.99"
$ws.Range("C2").Value = "$-14.11"
$ws.Range("D2").Value = "$-3.45"
$ws.Range("E2").Value = "$-4.99"
$ws.Range("F2").Value = "$0.00"
$ws.Range("G2").Value = "$0.00"
$ws.Range("H2").Value = "$0.00"
$ws.Range("I2").Value = "$0.00"
$ws.Range("J2").Value = "$-0.20"
$ws.Range("K2").Value = "$0.00"
$ws.Range("L2").Value = "$0.00"
$ws.Range("M2").Value = "$0.00"
$ws.Range("N2").Value = "1"
$ws.Range("O2").Value = "$0.24"
$ws.Range("P2").Value = "1.04%"

# Row 3: CONCORD
$ws.Range("A3").Value = "CONCORD"
$ws.Range("B3").Value = "$0.00"
$ws.Range("C3").Value = "$0.00"
$ws.Range("D3").Value = "$0.00"
$ws.Range("E3").Value = "$0.00"
$ws.Range("F3").Value = "$0.00"
$ws.Range("G3").Value = "$0.00"
$ws.Range("H3").Value = "$0.00"
$ws.Range("I3").Value = "$0.00"
$ws.Range("J3").Value = "$-0.99"
$ws.Range("K3").Value = "$0.00"
$ws.Range("L3").Value = "$0.00"
$ws.Range("M3").Value = "$0.00"
$ws.Range("N3").Value = "0"
$ws.Range("O3").Value = "$-0.99"
$ws.Range("P3").Value = "0.00%"

# Row 4: CRICUT
$ws.Range("A4").Value = "CRICUT"
$ws.Range("B4").Value = "$2,060.64"
$ws.Range("C4").Value = "$-441.25"
$ws.Range("D4").Value = "$-309.43"
$ws.Range("E4").Value = "$-740.32"
$ws.Range("F4").Value = "$0.00"
$ws.Range("G4").Value = "$-18.91"
$ws.Range("H4").Value = "$0.00"
$ws.Range("I4").Value = "$-0.55"
$ws.Range("J4").Value = "$-171.59"
$ws.Range("K4").Value = "$-0.18"
$ws.Range("L4").Value = "$-70.71"
$ws.Range("M4").Value = "$-17.68"
$ws.Range("N4").Value = "190"
$ws.Range("O4").Value = "$379.64"
$ws.Range("P4").Value = "18.42%"

# Row 5: HY-C
$ws.Range("A5").Value = "HY-C"
$ws.Range("B5").Value = "$94.99"
$ws.Range("C5").Value = "$-28.67"
$ws.Range("D5").Value = "$-14.25"
$ws.Range("E5").Value = "$-23.67"
$ws.Range("F5").Value = "$0.00"
$ws.Range("G5").Value = "$0.00"
$ws.Range("H5").Value = "$0.00"
$ws.Range("I5").Value = "$0.00"
$ws.Range("J5").Value = "$-593.50"
$ws.Range("K5").Value = "$-0.01"
$ws.Range("L5").Value = "$-152.32"
$ws.Range("M5").Value = "$-76.16"
$ws.Range("N5").Value = "1"
$ws.Range("O5").Value = "$-641.27"
$ws.Range("P5").Value = "-675.09%"

# Row 6: LITANIA
$ws.Range("A6").Value = "LITANIA"
$ws.Range("B6").Value = "$22,147.35"
$ws.Range("C6").Value = "$-11,708.80"
$ws.Range("D6").Value = "$-3,567.35"
$ws.Range("E6").Value = "$-529.72"
$ws.Range("F6").Value = "$0.00"
$ws.Range("G6").Value = "$-277.42"
$ws.Range("H6").Value = "$0.00"
$ws.Range("I6").Value = "$0.00"
$ws.Range("J6").Value = "$-90.33"
$ws.Range("K6").Value = "$-2.14"
$ws.Range("L6").Value = "$-216.63"
$ws.Range("M6").Value = "$-54.15"
$ws.Range("N6").Value = "85"
$ws.Range("O6").Value = "$6,194.86"
$ws.Range("P6").Value = "27.97%"

# Row 7: MWI-VETONE
$ws.Range("A7").Value = "MWI-VETONE"
$ws.Range("B7").Value = "$24.99"
$ws.Range("C7").Value = "$-13.90"
$ws.Range("D7").Value = "$-3.75"
$ws.Range("E7").Value = "$-3.68"
$ws.Range("F7").Value = "$0.00"
$ws.Range("G7").Value = "$0.00"
$ws.Range("H7").Value = "$0.00"
$ws.Range("I7").Value = "$0.00"
$ws.Range("J7").Value = "$0.00"
$ws.Range("K7").Value = "$0.00"
$ws.Range("L7").Value = "$0.00"
$ws.Range("M7").Value = "$0.00"
$ws.Range("N7").Value = "1"
$ws.Range("O7").Value = "$3.66"
$ws.Range("P7").Value = "14.65%"

# Row 8: PATRICIA NASH
$ws.Range("A8").Value = "PATRICIA NASH"
$ws.Range("B8").Value = "$7,216.00"
$ws.Range("C8").Value = "$-3,312.00"
$ws.Range("D8").Value = "$-1,062.72"
$ws.Range("E8").Value = "$-267.91"
$ws.Range("F8").Value = "$0.00"
$ws.Range("G8").Value = "$-59.47"
$ws.Range("H8").Value = "$0.00"
$ws.Range("I8").Value = "$-5.61"
$ws.Range("J8").Value = "$-117.40"
$ws.Range("K8").Value = "$-0.70"
$ws.Range("L8").Value = "$-317.24"
$ws.Range("M8").Value = "$-317.24"
$ws.Range("N8").Value = "46"
$ws.Range("O8").Value = "$2,132.42"
$ws.Range("P8").Value = "29.55%"

# Row 9: VIROX
$ws.Range("A9").Value = "VIROX"
$ws.Range("B9").Value = "$271,751.40"
$ws.Range("C9").Value = "$-98,771.67"
$ws.Range("D9").Value = "$-39,128.88"
$ws.Range("E9").Value = "$-55,460.44"
$ws.Range("F9").Value = "$0.00"
$ws.Range("G9").Value = "$-13,830.67"
$ws.Range("H9").Value = "$-13,830.67"
$ws.Range("I9").Value = "$-1,301.85"
$ws.Range("J9").Value = "$-1,346.00"
$ws.Range("K9").Value = "$-26.31"
$ws.Range("L9").Value = "$-3,773.79"
$ws.Range("M9").Value = "$-3,773.79"
$ws.Range("N9").Value = "1,053"
$ws.Range("O9").Value = "$58,111.79"
$ws.Range("P9").Value = "21.38%"

# Row 10: WABASH VALLEY FARMS
$ws.Range("A10").Value = "WABASH VALLEY FARMS"
$ws.Range("B10").Value = "$109,799.32"
$ws.Range("C10").Value = "$-44,673.16"
$ws.Range("D10").Value = "$-16,366.10"
$ws.Range("E10").Value = "$-20,540.63"
$ws.Range("F10").Value = "$0.00"
$ws.Range("G10").Value = "$-2,301.96"
$ws.Range("H10").Value = "$-2,301.96"
$ws.Range("I10").Value = "$0.00"
$ws.Range("J10").Value = "$-3,034.65"
$ws.Range("K10").Value = "$-10.59"
$ws.Range("L10").Value = "$-1,994.42"
$ws.Range("M10").Value = "$-997.20"
$ws.Range("N10").Value = "2,666"
$ws.Range("O10").Value = "$21,875.03"
$ws.Range("P10").Value = "19.92%"
